# Aanvullingen module4 - figuren en backup aia files
#
# 1) Update the cached "datetimeFigureOut" footer-date field (shown as the
#    literal text "12/10/2018") to "1/11/2019" on the Slide Master and on
#    every Custom Layout's "Date Placeholder" shape.
# 2) Split the single run "Hartslag monitor altijd op zak" on slide 1 into
#    two runs with identical formatting: "Je hartslagmonitor " + "altijd op zak".

$p = $ppt.ActivePresentation
$newDate = "1/11/2019"

function Update-DatePlaceholder {
    param($shapes)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster

# Slide master footer date placeholder.
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout has its own footer date placeholder.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1: "Hartslag monitor altijd op zak" -> two runs:
# "Je hartslagmonitor " + "altijd op zak" (same run formatting throughout).
$slide = $p.Slides.Item(1)
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $shp = $slide.Shapes.Item($k)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Hartslag monitor altijd op zak") {
            $prefixLen = "Hartslag monitor ".Length
            $prefixRange = $tr.Characters(1, $prefixLen)
            $prefixRange.Text = "Je hartslagmonitor "
        }
    }
}
